$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.841834217421933
$ws.Range("C2").Value = 0.2377829919085457
$ws.Range("D2").Value = 0.1132707645347821
$ws.Range("E2").Value = 0.0490461204128616
$ws.Range("F2").Value = 2.421976638479862
$ws.Range("L2").Value = 0.2716327507868357
$ws.Range("N2").Value = 1.946590460329148
$ws.Range("B3").Value = 2.6830833957452
$ws.Range("C3").Value = 0.2069324353719253
$ws.Range("D3").Value = 0.113647159158262
$ws.Range("E3").Value = 0.04924918981995363
$ws.Range("F3").Value = 2.364479096077233
$ws.Range("L3").Value = 0.2616684373782334
$ws.Range("N3").Value = 1.95822297313741
$ws.Range("B4").Value = 2.58746215276841
$ws.Range("C4").Value = 0.1880198126431765
$ws.Range("D4").Value = 0.1139298444665897
$ws.Range("E4").Value = 0.04938467912766198
$ws.Range("F4").Value = 2.330742348702174
$ws.Range("L4").Value = 0.2557277015415167
$ws.Range("N4").Value = 1.966024942044008
$ws.Range("B5").Value = 2.548957376177611
$ws.Range("C5").Value = 0.1803193773976091
$ws.Range("D5").Value = 0.1140579358833165
$ws.Range("E5").Value = 0.04944261196031796
$ws.Range("F5").Value = 2.317385159484473
$ws.Range("L5").Value = 0.2533511093697314
$ws.Range("N5").Value = 1.969369320994161
$ws.Range("B6").Value = 2.542591450739849
$ws.Range("C6").Value = 0.1790411002433245
$ws.Range("D6").Value = 0.1140799820172802
$ws.Range("E6").Value = 0.04945239602791851
$ws.Range("F6").Value = 2.315190718983203
$ws.Range("L6").Value = 0.252959145343425
$ws.Range("N6").Value = 1.969934594587052
$ws.Range("B7").Value = 2.586941000067554
$ws.Range("C7").Value = 0.1879159362466964
$ws.Range("D7").Value = 0.1139315198388005
$ws.Range("E7").Value = 0.04938544941335099
$ws.Range("F7").Value = 2.330560630820401
$ws.Range("L7").Value = 0.2556954709874617
$ws.Range("N7").Value = 1.966069378471573
$ws.Range("B8").Value = 2.786709895398246
$ws.Range("C8").Value = 0.2271390126917083
$ws.Range("D8").Value = 0.113389795120149
$ws.Range("E8").Value = 0.04911389825117363
$ws.Range("F8").Value = 2.40182428752982
$ws.Range("L8").Value = 0.2681600511346147
$ws.Range("N8").Value = 1.950463952926356
$ws.Range("B9").Value = 3.193367121953315
$ws.Range("C9").Value = 0.3043344642938735
$ws.Range("D9").Value = 0.1127401793202409
$ws.Range("E9").Value = 0.04866697884752202
$ws.Range("F9").Value = 2.554161799494267
$ws.Range("L9").Value = 0.2940259638491085
$ws.Range("N9").Value = 1.925133008421852
$ws.Range("B10").Value = 3.501561083297929
$ws.Range("C10").Value = 0.3612856057851559
$ws.Range("D10").Value = 0.1125194272343464
$ws.Range("E10").Value = 0.04839063558225565
$ws.Range("F10").Value = 2.673992092429245
$ws.Range("L10").Value = 0.3139209947143939
$ws.Range("N10").Value = 1.909788611674941
$ws.Range("B11").Value = 3.643888120940915
$ws.Range("C11").Value = 0.3872600052351913
$ws.Range("D11").Value = 0.1124757954355928
$ws.Range("E11").Value = 0.04827617917387617
$ws.Range("F11").Value = 2.730275165509028
$ws.Range("L11").Value = 0.3231707290507444
$ws.Range("N11").Value = 1.903528881311615
$ws.Range("B12").Value = 3.698095256445754
$ws.Range("C12").Value = 0.3971065752492109
$ws.Range("D12").Value = 0.1124675299822329
$ws.Range("E12").Value = 0.04823445324865983
$ws.Range("F12").Value = 2.751846892124007
$ws.Range("L12").Value = 0.3267024377778966
$ws.Range("N12").Value = 1.901263075431373
$ws.Range("B13").Value = 3.686406868919789
$ws.Range("C13").Value = 0.3949854515708466
$ws.Range("D13").Value = 0.1124689414022484
$ws.Range("E13").Value = 0.04824336782185501
$ws.Range("F13").Value = 2.747189474832965
$ws.Range("L13").Value = 0.3259405238873541
$ws.Range("N13").Value = 1.901746388098829
$ws.Range("B14").Value = 3.648341513346281
$ws.Range("C14").Value = 0.3880698697329876
$ws.Range("D14").Value = 0.1124749494249926
$ws.Range("E14").Value = 0.04827271397645738
$ws.Range("F14").Value = 2.732044680831649
$ws.Range("L14").Value = 0.3234607008326833
$ws.Range("N14").Value = 1.903340369560411
$ws.Range("B15").Value = 3.625066031417532
$ws.Range("C15").Value = 0.3838352899972506
$ws.Range("D15").Value = 0.1124797075330548
$ws.Range("E15").Value = 0.04829089976662093
$ws.Range("F15").Value = 2.722801841294597
$ws.Range("L15").Value = 0.3219455297760447
$ws.Range("N15").Value = 1.90433038381741
$ws.Range("B16").Value = 3.492302920214797
$ws.Range("C16").Value = 0.3595895507752971
$ws.Range("D16").Value = 0.1125234291382924
$ws.Range("E16").Value = 0.04839834185325831
$ws.Range("F16").Value = 2.670349833742449
$ws.Range("L16").Value = 0.3133205458589572
$ws.Range("N16").Value = 1.910212284256005
$ws.Range("B17").Value = 3.41140559272884
$ws.Range("C17").Value = 0.3447334429769739
$ws.Range("D17").Value = 0.1125648600622284
$ws.Range("E17").Value = 0.04846713492843602
$ws.Range("F17").Value = 2.638628664890859
$ws.Range("L17").Value = 0.3080807210010477
$ws.Range("N17").Value = 1.914005893466907
$ws.Range("B18").Value = 3.365075588673619
$ws.Range("C18").Value = 0.3361948525409844
$ws.Range("D18").Value = 0.1125940297990908
$ws.Range("E18").Value = 0.04850776223513265
$ws.Range("F18").Value = 2.620550015391245
$ws.Range("L18").Value = 0.3050856571226888
$ws.Range("N18").Value = 1.916255633949376
$ws.Range("B19").Value = 3.349423246077038
$ws.Range("C19").Value = 0.3333048761396071
$ws.Range("D19").Value = 0.1126048205428773
$ws.Range("E19").Value = 0.04852169995038835
$ws.Range("F19").Value = 2.614457382392146
$ws.Range("L19").Value = 0.3040747884484176
$ws.Range("N19").Value = 1.917028964194671
$ws.Range("B20").Value = 3.419996516626099
$ws.Range("C20").Value = 0.3463142470774301
$ws.Range("D20").Value = 0.1125598963739378
$ws.Range("E20").Value = 0.04845970216437312
$ws.Range("F20").Value = 2.64198817527739
$ws.Range("L20").Value = 0.3086365666803061
$ws.Range("N20").Value = 1.913595037459899
$ws.Range("B21").Value = 3.659513751075565
$ws.Range("C21").Value = 0.3901008473268348
$ws.Range("D21").Value = 0.1124729598998826
$ws.Range("E21").Value = 0.0482640504507712
$ws.Range("F21").Value = 2.736486026520964
$ws.Range("L21").Value = 0.324188294030634
$ws.Range("N21").Value = 1.902869331136145
$ws.Range("B22").Value = 3.817867099816567
$ws.Range("C22").Value = 0.4187803459709585
$ws.Range("D22").Value = 0.1124643131138967
$ws.Range("E22").Value = 0.04814560052720473
$ws.Range("F22").Value = 2.799754657441071
$ws.Range("L22").Value = 0.3345216605272299
$ws.Range("N22").Value = 1.896469736218322
$ws.Range("B23").Value = 3.733183194872481
$ws.Range("C23").Value = 0.4034675144939115
$ws.Range("D23").Value = 0.1124644899707832
$ws.Range("E23").Value = 0.0482079582073851
$ws.Range("F23").Value = 2.765847672990787
$ws.Range("L23").Value = 0.3289909245547022
$ws.Range("N23").Value = 1.899829142756928
$ws.Range("B24").Value = 3.416112002924649
$ws.Range("C24").Value = 0.3455995581924753
$ws.Range("D24").Value = 0.1125621237988383
$ws.Range("E24").Value = 0.04846305916149518
$ws.Range("F24").Value = 2.640468847931771
$ws.Range("L24").Value = 0.3083852147376547
$ws.Range("N24").Value = 1.913780571414492
$ws.Range("B25").Value = 3.081725240779633
$ws.Range("C25").Value = 0.2834144603457389
$ws.Range("D25").Value = 0.1128712313765021
$ws.Range("E25").Value = 0.04877873720236714
$ws.Range("F25").Value = 2.511578800826129
$ws.Range("L25").Value = 0.2868737749715677
$ws.Range("N25").Value = 1.931415975659036
